$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update B4: output location path
$ws.Range("B4").Value = "C:\Users\dpere\Documents\JTMT\Projects\תחזיות_דמוגרפיות\קבצי עבודה\142_מתחם_אנגל\בהת\For_approval\Reference_tabels"

# Update B5: boolean flag from FALSE to TRUE
$ws.Range("B5").Value = $true

# Update B6: new layer location path
$ws.Range("B6").Value = "C:\Users\dpere\Documents\JTMT\Projects\תחזיות_דמוגרפיות\קבצי עבודה\142_מתחם_אנגל\בהת\For_approval\Reference_tabels\shp\TAZ_V4_240627_with_geo_info.shp"
